$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = "4.94-grazing"
    4  = "0.6587-grazing, scrub"
    5  = "3.4992-grazing, foiliage, scrub"
    6  = "3.2933-grazing"
    7  = "0.5763-grazing, foiliage, scrub, seeds"
    8  = "0.494-meat"
    9  = "0.6175-grazing, scrub, foilage"
    10 = "5.7633-grazing"
    11 = "32.9333-grazing, foiliage, scrub"
    12 = "54.34-grazing, foiliage, scrub"
    13 = "0.4528-grazing"
    14 = "2.0995-meat"
    15 = "0.2058-meat"
    16 = "0.0463-meat"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
